$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 953.375
$ws.Range("I2").Value = 604.5
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 604.5
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -491.5
$ws.Range("N2").Value = -2226
$ws.Range("H6").Value = 47.357143
$ws.Range("I6").Value = 38.583332
$ws.Range("K6").Value = 115.749996
$ws.Range("M6").Value = -3.749995999999996
$ws.Range("H40").Value = 2534.75
$ws.Range("J40").Value = 2948.4285
$ws.Range("L40").Value = 2948.4285
$ws.Range("N40").Value = -3298.4285
$ws.Range("H125").Value = 7166.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 7166.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 64498.5
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -69418.5
$ws.Range("H132").Value = 5095.625
$ws.Range("I132").Value = 5258.048
$ws.Range("J132").Value = 3958.6667
$ws.Range("K132").Value = 15774.144
$ws.Range("L132").Value = 11876.0001
$ws.Range("M132").Value = -13244.144
$ws.Range("N132").Value = -16936.0001
$ws.Range("H138").Value = 6807525
$ws.Range("J138").Value = 10759014
$ws.Range("L138").Value = 32277042
$ws.Range("N138").Value = -32287322
$ws.Range("H141").Value = 3170.8
$ws.Range("I141").Value = 3170.8
$ws.Range("K141").Value = 9512.400000000001
$ws.Range("M141").Value = -4332.400000000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 333333540
$ws.Range("I8").Value = 500000000
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 500000000
$ws.Range("L8").Value = 600
$ws.Range("M8").Value = -499999856
$ws.Range("N8").Value = -888
$ws.Range("H63").Value = 5826
$ws.Range("I63").Value = 1652.1666
$ws.Range("K63").Value = 1652.1666
$ws.Range("M63").Value = -966.1666
$ws.Range("H66").Value = 5826
$ws.Range("I66").Value = 1652.1666
$ws.Range("K66").Value = 8260.833000000001
$ws.Range("M66").Value = -4828.833000000001
$ws.Range("H124").Value = 40036.168
$ws.Range("J124").Value = 40036.168
$ws.Range("L124").Value = 40036.168
$ws.Range("N124").Value = -49856.168
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3167.8572
$ws.Range("I20").Value = 4491.6665
$ws.Range("K20").Value = 4491.6665
$ws.Range("M20").Value = -4244.6665
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4997.8335
$ws.Range("I62").Value = 3624.25
$ws.Range("J62").Value = 7745
$ws.Range("K62").Value = 3624.25
$ws.Range("L62").Value = 7745
$ws.Range("M62").Value = -3000.25
$ws.Range("N62").Value = -8993
$ws.Range("H65").Value = 4997.8335
$ws.Range("I65").Value = 3624.25
$ws.Range("J65").Value = 7745
$ws.Range("K65").Value = 18121.25
$ws.Range("L65").Value = 38725
$ws.Range("M65").Value = -15001.25
$ws.Range("N65").Value = -44965
$ws.Range("H105").Value = 8080
$ws.Range("I105").Value = 2098.5
$ws.Range("J105").Value = 13396.889
$ws.Range("K105").Value = 2098.5
$ws.Range("L105").Value = 13396.889
$ws.Range("M105").Value = -351.5
$ws.Range("N105").Value = -16890.889
$ws.Range("H132").Value = 78156.82000000001
$ws.Range("I132").Value = 124674.94
$ws.Range("J132").Value = 6265.1816
$ws.Range("K132").Value = 374024.82
$ws.Range("L132").Value = 18795.5448
$ws.Range("M132").Value = -371494.82
$ws.Range("N132").Value = -23855.5448
$ws.Range("H134").Value = 1463.1364
$ws.Range("I134").Value = 1076.5555
$ws.Range("K134").Value = 3229.6665
$ws.Range("M134").Value = -694.6664999999998
$ws.Range("H141").Value = 288854.34
$ws.Range("J141").Value = 309961.12
$ws.Range("L141").Value = 309961.12
$ws.Range("N141").Value = -320321.12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3180.625
$ws.Range("I3").Value = 896.5
$ws.Range("J3").Value = 10033
$ws.Range("K3").Value = 2689.5
$ws.Range("L3").Value = 30099
$ws.Range("M3").Value = -2577.5
$ws.Range("N3").Value = -30323
$ws.Range("H4").Value = 35873830
$ws.Range("I4").Value = 59130736
$ws.Range("J4").Value = 16845454
$ws.Range("K4").Value = 177392208
$ws.Range("L4").Value = 50536362
$ws.Range("M4").Value = -177392096
$ws.Range("N4").Value = -50536586
$ws.Range("H13").Value = 1668483
$ws.Range("I13").Value = 5000004.5
$ws.Range("J13").Value = 2722.25
$ws.Range("K13").Value = 15000013.5
$ws.Range("L13").Value = 8166.75
$ws.Range("M13").Value = -14999845.5
$ws.Range("N13").Value = -8502.75
$ws.Range("H68").Value = 1298.3334
$ws.Range("I68").Value = 1096.6666
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 3289.9998
$ws.Range("L68").Value = 4500
$ws.Range("M68").Value = -2478.9998
$ws.Range("N68").Value = -6122
$ws.Range("H71").Value = 1298.3334
$ws.Range("I71").Value = 1096.6666
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 9869.999400000001
$ws.Range("L71").Value = 13500
$ws.Range("M71").Value = -5813.999400000001
$ws.Range("N71").Value = -21612
$ws.Range("H131").Value = 1475.3889
$ws.Range("J131").Value = 1849.55
$ws.Range("L131").Value = 5548.65
$ws.Range("N131").Value = -15628.65
$ws.Range("H133").Value = 11553.9
$ws.Range("I133").Value = 7943.143
$ws.Range("K133").Value = 23829.429
$ws.Range("M133").Value = -18769.429
$ws.Range("H134").Value = 3700.4666
$ws.Range("I134").Value = 2536.2144
$ws.Range("K134").Value = 7608.6432
$ws.Range("M134").Value = -2538.6432
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 925
$ws.Range("I9").Value = 925
$ws.Range("K9").Value = 925
$ws.Range("M9").Value = -755
$ws.Range("H70").Value = 506754
$ws.Range("I70").Value = 672672
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 672672
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -672402
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 506754
$ws.Range("I73").Value = 672672
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 672672
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -671736
$ws.Range("N73").Value = -10872
$ws.Range("H132").Value = 4384.6665
$ws.Range("I132").Value = 3762.6
$ws.Range("K132").Value = 11287.8
$ws.Range("M132").Value = -8757.799999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 60000
$ws.Range("I63").Value = 60000
$ws.Range("K63").Value = 60000
$ws.Range("M63").Value = -59251
$ws.Range("H66").Value = 60000
$ws.Range("I66").Value = 60000
$ws.Range("K66").Value = 180000
$ws.Range("M66").Value = -176256
$ws.Range("H68").Value = 3163.1614
$ws.Range("I68").Value = 2789.423
$ws.Range("J68").Value = 5106.6
$ws.Range("K68").Value = 2789.423
$ws.Range("L68").Value = 5106.6
$ws.Range("M68").Value = -2040.423
$ws.Range("N68").Value = -6604.6
$ws.Range("H71").Value = 3163.1614
$ws.Range("I71").Value = 2789.423
$ws.Range("J71").Value = 5106.6
$ws.Range("K71").Value = 13947.115
$ws.Range("L71").Value = 25533
$ws.Range("M71").Value = -10203.115
$ws.Range("N71").Value = -33021
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 4135.4595
$ws.Range("I122").Value = 3771.3333
$ws.Range("J122").Value = 4807.6924
$ws.Range("K122").Value = 11313.9999
$ws.Range("L122").Value = 14423.0772
$ws.Range("M122").Value = -8863.999899999999
$ws.Range("N122").Value = -19323.0772
$ws.Range("H127").Value = 67900
$ws.Range("J127").Value = 67900
$ws.Range("L127").Value = 67900
$ws.Range("N127").Value = -77820
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H70").Value = 49199.8
$ws.Range("I70").Value = 40249.75
$ws.Range("K70").Value = 40249.75
$ws.Range("M70").Value = -39934.75
$ws.Range("H73").Value = 49199.8
$ws.Range("I73").Value = 40249.75
$ws.Range("K73").Value = 40249.75
$ws.Range("M73").Value = -39157.75
$ws.Range("H113").Value = 508.56522
$ws.Range("I113").Value = 259.41177
$ws.Range("J113").Value = 1214.5
$ws.Range("K113").Value = 778.23531
$ws.Range("L113").Value = 3643.5
$ws.Range("M113").Value = 1391.76469
$ws.Range("N113").Value = -7983.5
$ws.Range("H132").Value = 4123.1587
$ws.Range("J132").Value = 6997.5
$ws.Range("L132").Value = 20992.5
$ws.Range("N132").Value = -26052.5
$ws.Range("H136").Value = 1864.7188
$ws.Range("I136").Value = 1842.5
$ws.Range("J136").Value = 1961
$ws.Range("K136").Value = 5527.5
$ws.Range("L136").Value = 5883
$ws.Range("M136").Value = -2977.5
$ws.Range("N136").Value = -10983
